$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a rolling daily price log: row 2 is always "today", and each
# day a new row is inserted at the top pushing every existing row down by
# one, with the oldest row (106) falling out the bottom by being duplicated
# into the new row 107 (no separate "drop" step is needed for that - it
# falls straight out of the shift).

$lastRow = 106
$newLastRow = $lastRow + 1

# 1) Clone row 106's cell formatting (number formats, borders, alignment)
#    into the brand-new row 107 so the appended row matches the table style.
$ws.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$ws.Range("A" + $newLastRow + ":F" + $newLastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Columns A (Date) and E (Circular Date) hold dd-mm-yyyy text. Force them
#    to Text format before the bulk write below so values like "01-11-2025"
#    aren't auto-parsed as a US-style m-d-yyyy date by the Value2 setter.
$ws.Range("A3:A" + $newLastRow).NumberFormat = "@"
$ws.Range("E3:E" + $newLastRow).NumberFormat = "@"

# 3) Shift every existing data row (2..106) down by one row (3..107) in one
#    bulk copy, preserving each row's Date/Description/Code/Price/Circular
#    Date/Link together.
$arr = $ws.Range("A2:F" + $lastRow).Value2
$ws.Range("A3:F" + $newLastRow).Value2 = $arr

# 4) Re-apply the normal (General-format) cell style to columns A and E so
#    they match the rest of the table instead of staying tagged as Text.
$ws.Range("C3:C" + $newLastRow).Copy()
$ws.Range("A3:A" + $newLastRow).PasteSpecial(-4122)
$ws.Range("E3:E" + $newLastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5) Insert today's new row 2. Only the Date advances by one day; the rest
#    of the circular/price info is unchanged until Nalco issues a new one.
$prevDateText = $ws.Cells.Item(3, 1).Value2
$p = $prevDateText.Split("-")
$prevDate = Get-Date -Year ([int]$p[2]) -Month ([int]$p[1]) -Day ([int]$p[0])
$newDate = $prevDate.AddDays(1)
$newDateText = $newDate.ToString("dd-MM-yyyy")

$ws.Cells.Item(2, 1).Value = $newDateText
$ws.Cells.Item(2, 2).Value = $ws.Cells.Item(3, 2).Value2
$ws.Cells.Item(2, 3).Value = $ws.Cells.Item(3, 3).Value2
$ws.Cells.Item(2, 4).Value = $ws.Cells.Item(3, 4).Value2
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(3, 5).Value2
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(3, 6).Value2

# 6) The hyperlink collection doesn't move with the rows it's attached to,
#    so rebuild it from scratch against the new layout: wipe everything
#    (deleting any one hyperlink clears the whole sheet's collection in this
#    host) then re-add one hyperlink per data row, F2 down to F107, each
#    pointing at that row's own circular link text.
$ws.Range("F2").Hyperlinks.Delete()
for ($r = 2; $r -le $newLastRow; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Range("F" + $r), $target) | Out-Null
}

Write-Output "done"
